$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for years 2000-2009 (old rows 2-11); this shifts the
# 2010-2020 data up so it now starts at row 2.
$ws.Rows("2:11").Delete()

# After the delete, rows 2-12 hold years 2010-2020. A handful of values
# for 2015, 2016, 2018, 2019 and 2020 were re-rounded to whole numbers.
$ws.Range("B7").Value = 472556
$ws.Range("C7").Value = 80699

$ws.Range("B8").Value = 493254
$ws.Range("C8").Value = 78307

$ws.Range("B10").Value = 578244
$ws.Range("C10").Value = 92322

$ws.Range("B11").Value = 550530
$ws.Range("C11").Value = 100943

$ws.Range("B12").Value = 566181
$ws.Range("C12").Value = 103471

# Append the new 2021 row. Copy the formatting of the row above (bold,
# centered/top aligned, thin border) onto the new label cell first, then
# fill in the values.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "2021年"

$ws.Range("B13").Value = 593225.61
$ws.Range("C13").Value = 118783.71
